$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) values
$ws.Range("B2").Value = 7.8053612560655665
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.1180024452172574
$ws.Range("E2").Value = 1.5280757406754049

# Row 3 (STR) values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 12.864998719718898
$ws.Range("D3").Value = 11.223432451595187
$ws.Range("E3").Value = -3.2063479174769896

# Update selection to match new active range
$ws.Range("B1:E3").Select()
